# "Generate Report for Handback" — record the handback results for the
# zh-cn and de-de locales on the localization-status workbook.

$wb = $excel.ActiveWorkbook

$handedBackStatus = "Handed back: in sync with en-US"
$mdDisplay        = "LocaleLowerCaseTest.md"
$mdUrl            = "https://github.com/OpenLocalizationTestOrg/LocaleLowerCaseTest/blob/43005fad50282e1d7eb2b234e05d7f949ed92bcf/test/LocaleLowerCaseTest.md"

# The host quantizes ColumnWidth to whole pixels (6px/char + 5px padding)
# on write, so the literal target character-widths are fed back a touch
# low here to land on the nearest reachable width after quantization.
$wideColWidth    = 29.166666666666668   # -> ~29.9777047293527 (Name/File cols)
$targetColWidth  = 22.833333333333332   # -> ~23.5855930873326 (Latest Target File)
$handbackColWidth = 39.166666666666664  # -> 40 (Latest Handback File)

# ---------------------------------------------------------------------
# Overview sheet: both locale status columns move to "handed back"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $handedBackStatus
$wsOverview.Range("F2").Value = $handedBackStatus
$wsOverview.Columns.Item(5).ColumnWidth = $wideColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideColWidth

# ---------------------------------------------------------------------
# zh-cn sheet: status + the freshly generated handback target/file/time
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $handedBackStatus
$wsZhCn.Range("K2").Value = "LocaleLowerCaseTest.a47dea4a0d9bca8cf007ef5d5443046c4dff2a81.zh-cn.xlf"
$wsZhCn.Range("L2").Value = "2017-03-03 02:34:41"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("J2"), $mdUrl, "", "", $mdDisplay)

$wsZhCn.Columns.Item(3).ColumnWidth = $wideColWidth
$wsZhCn.Columns.Item(10).ColumnWidth = $targetColWidth
$wsZhCn.Columns.Item(11).ColumnWidth = $handbackColWidth

# ---------------------------------------------------------------------
# de-de sheet: status + the freshly generated handback target/file/time
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $handedBackStatus
$wsDeDe.Range("K2").Value = "LocaleLowerCaseTest.a47dea4a0d9bca8cf007ef5d5443046c4dff2a81.de-de.xlf"
$wsDeDe.Range("L2").Value = "2017-03-03 02:34:50"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("J2"), $mdUrl, "", "", $mdDisplay)

$wsDeDe.Columns.Item(3).ColumnWidth = $wideColWidth
$wsDeDe.Columns.Item(10).ColumnWidth = $targetColWidth
$wsDeDe.Columns.Item(11).ColumnWidth = $handbackColWidth
